$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old data rows (2-5); header row (1) stays untouched.
$ws.Range("A2:Q5").ClearContents()

# --- Row 2: FEVar / FEATV ---
$ws.Range("A2").Value = "FEVar"
$ws.Range("B2").Value = "FEATV"
$ws.Range("F2").Value = 0.09
$ws.Range("G2").Value = 2.77
$ws.Range("H2").Value = 0.093
$ws.Range("I2").Value = 1.408
$ws.Range("J2").Value = 0.911
$ws.Range("K2").Value = 0.422
$ws.Range("L2").NumberFormat = "0.0"
$ws.Range("L2").Value = 3.4
$ws.Range("M2").NumberFormat = "0.0"
$ws.Range("M2").Value = 15.39
$ws.Range("N2").NumberFormat = "0.000"
$ws.Range("N2").Value = 3.397487
$ws.Range("O2").NumberFormat = "0.000"
$ws.Range("O2").Value = 15.394725
$ws.Range("P2").NumberFormat = "0.000"
$ws.Range("P2").Value = 0.99672
$ws.Range("Q2").NumberFormat = "0.000"
$ws.Range("Q2").Value = 0.026543

# --- Row 3: FEVar / FEATV / DisgVar / DisgATV ---
$ws.Range("A3").Value = "FEVar"
$ws.Range("B3").Value = "FEATV"
$ws.Range("C3").Value = "DisgVar"
$ws.Range("D3").Value = "DisgATV"
$ws.Range("F3").Value = 0.09
$ws.Range("G3").Value = 2.77
$ws.Range("H3").Value = 0.093
$ws.Range("I3").Value = 1.408
$ws.Range("J3").Value = 0.911
$ws.Range("K3").Value = 0.422
$ws.Range("L3").NumberFormat = "0.0"
$ws.Range("L3").Value = 3.4
$ws.Range("M3").NumberFormat = "0.0"
$ws.Range("M3").Value = 15.39
$ws.Range("N3").NumberFormat = "0.000"
$ws.Range("N3").Value = 3.397486
$ws.Range("O3").NumberFormat = "0.000"
$ws.Range("O3").Value = 15.394713
$ws.Range("P3").NumberFormat = "0.000"
$ws.Range("P3").Value = 0.99672
$ws.Range("Q3").NumberFormat = "0.000"
$ws.Range("Q3").Value = 0.026543

# --- Row 4: FEVar / FEATV / DisgVar / DisgATV / Var ---
$ws.Range("A4").Value = "FEVar"
$ws.Range("B4").Value = "FEATV"
$ws.Range("C4").Value = "DisgVar"
$ws.Range("D4").Value = "DisgATV"
$ws.Range("E4").Value = "Var"
$ws.Range("F4").Value = 0.14
$ws.Range("G4").Value = 3.85
$ws.Range("H4").Value = 0.133
$ws.Range("I4").Value = 1.359
$ws.Range("J4").Value = 0.911
$ws.Range("K4").Value = 0.422
$ws.Range("L4").NumberFormat = "0.0"
$ws.Range("L4").Value = 4.86
$ws.Range("M4").NumberFormat = "0.0"
$ws.Range("M4").Value = 22.37
$ws.Range("N4").NumberFormat = "0.000"
$ws.Range("N4").Value = 4.859968
$ws.Range("O4").NumberFormat = "0.000"
$ws.Range("O4").Value = 22.36728
$ws.Range("P4").NumberFormat = "0.000"
$ws.Range("P4").Value = 0.996721
$ws.Range("Q4").NumberFormat = "0.000"
$ws.Range("Q4").Value = 0.026542

# Selection / view now covers the (now smaller) used range A1:Q4.
$ws.Range("A1:Q4").Select()

Write-Host "Rewrote NI_Est data rows 2-4 with updated estimates."
